$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "Is Active" column (F2:F5) from text "True"/"False" shared
# strings into native Excel Boolean values.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
